# Update the workbook "Översikt ESLÖV":
#  1) Change the "Förändrad" (changed) date in column C for every existing
#     data row (rows 2-76) from 2023-09-08 (45177) to 2023-09-09 (45178).
#  2) Touch the row height of row 76 so it gets an explicit customHeight
#     flag like all the other data rows.
#  3) Append a brand-new data row (row 77) for case "A 42050-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) bump the "Förändrad" date for the existing rows (2..76) ---------
for ($r = 2; $r -le 76; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}

# --- 2) give row 76 the same explicit row height as its neighbours ------
$ws.Rows.Item(76).RowHeight = 15

# --- 3) append the new row 77 --------------------------------------------
$row = 77

$ws.Cells.Item($row, 1).Value = "A 42050-2023"

$ws.Cells.Item($row, 2).Value = 45177
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value = 45178
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value = "SKÅNE LÄN"
$ws.Cells.Item($row, 5).Value = "ESLÖV"

$ws.Cells.Item($row, 7).Value = 10.3
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

# Column R keeps the wrapped-text style used throughout the sheet, with no
# content (same as rows 75/76).
$ws.Cells.Item($row, 18).WrapText = $true
